$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'achilles compression'
$ws.Range("A2").Value = 'achilles support with compression'
$ws.Range("A3").Value = 'achilles wrap'
$ws.Range("A4").Value = 'achy joints'
$ws.Range("A5").Value = 'acls pocket'
$ws.Range("A6").Value = 'acting pro pants'
$ws.Range("A7").Value = 'acting pro women'
$ws.Range("A8").Value = 'active capri'
$ws.Range("A9").Value = 'active clothes'
$ws.Range("A10").Value = 'active clothes for women'
$ws.Range("A11").Value = 'active clothing'
$ws.Range("A12").Value = 'active clothing for women'
$ws.Range("A13").Value = 'active club leggings'
$ws.Range("A14").Value = 'active leggings for women'
$ws.Range("A15").Value = 'active leggings high waist'
$ws.Range("A16").Value = 'active pants'
$ws.Range("A17").Value = 'active plus size'
$ws.Range("A18").Value = 'active research compression pants'
$ws.Range("A19").Value = 'active research leggings'
$ws.Range("A20").Value = 'active research womens compression pants'
$ws.Range("A21").Value = 'active tights'
$ws.Range("A22").Value = 'active tights for women'
$ws.Range("A23").Value = 'active tights women'
$ws.Range("A24").Value = 'active wear'
$ws.Range("A25").Value = 'active wear capris for women'
$ws.Range("A26").Value = 'active wear clothes for women'
$ws.Range("A27").Value = 'active wear leggings for women'
$ws.Range("A28").Value = 'active wear leggings women'
$ws.Range("A29").Value = 'active wear tights'
$ws.Range("A30").Value = 'active wear woman'
$ws.Range("A31").Value = 'active wear women'
$ws.Range("A32").Value = 'active wear womens'
$ws.Range("A33").Value = 'active women wear'
$ws.Range("A34").Value = 'active womens clothing'
$ws.Range("A35").Value = 'adjustable dress form xl'
$ws.Range("A36").Value = 'adjustable hem tape'
$ws.Range("A37").Value = 'adult clothes for women'
$ws.Range("A38").Value = 'adult tights'
$ws.Range("A39").Value = 'aimado leggings for women'
$ws.Range("A40").Value = 'air compression for legs'
$ws.Range("A41").Value = 'air compression leg'
$ws.Range("A42").Value = 'air runner'
$ws.Range("A43").Value = 'airplane clothes for women'
$ws.Range("A44").Value = 'airplane clothing'
$ws.Range("A45").Value = 'airplane leggings'
$ws.Range("A46").Value = 'airplane leggings for women'
$ws.Range("A47").Value = 'airplane life jacket'
$ws.Range("A48").Value = 'airplane pants'
$ws.Range("A49").Value = 'airplane sleep support'
$ws.Range("A50").Value = 'airy pants'
$ws.Range("A51").Value = 'ajisai womens joggers pants'
$ws.Range("A52").Value = 'align apparel'
$ws.Range("A53").Value = 'align band'
$ws.Range("A54").Value = 'align crop'
$ws.Range("A55").Value = 'align leggings'
$ws.Range("A56").Value = 'align leggings women'
$ws.Range("A57").Value = 'align pants'
$ws.Range("A58").Value = 'align shorts women'
$ws.Range("A59").Value = 'all access leggings'
$ws.Range("A60").Value = 'all blacks jacket rugby'
$ws.Range("A61").Value = 'all blacks rugby women'
$ws.Range("A62").Value = 'all in motion leggings'
$ws.Range("A63").Value = 'all joint'
$ws.Range("A64").Value = 'all pro weight adjustable ankle weights'
$ws.Range("A65").Value = 'all weather generator'
$ws.Range("A66").Value = 'all weather jacket for women'
$ws.Range("A67").Value = 'all weather pants'
$ws.Range("A68").Value = 'all weather work pants for men'
$ws.Range("A69").Value = 'alternative apparel women pants'
$ws.Range("A70").Value = 'alternative apparel women shorts'
$ws.Range("A71").Value = 'althletic wear for women'
$ws.Range("A72").Value = 'always capri leggings'
$ws.Range("A73").Value = 'always leggings plus'
$ws.Range("A74").Value = 'always women leggings'
$ws.Range("A75").Value = 'amazing muscle'
$ws.Range("A76").Value = 'amazon cycling clothes'
$ws.Range("A77").Value = 'amazon international shopping'
$ws.Range("A78").Value = 'amazon knee braces'
$ws.Range("A79").Value = 'amazon leggings women'
$ws.Range("A80").Value = 'amazon logo tape'
$ws.Range("A81").Value = 'amazon pnws'
$ws.Range("A82").Value = 'amazon tights for women'
$ws.Range("A83").Value = 'amazon womens workout clothes'
$ws.Range("A84").Value = 'amazons choice leggings'
$ws.Range("A85").Value = 'amp wore'
$ws.Range("A86").Value = 'and still i rise'
$ws.Range("A87").Value = 'and support'
$ws.Range("A88").Value = 'andrew marc womens pants'
$ws.Range("A89").Value = 'angry stick man 5'
$ws.Range("A90").Value = 'angry woman'
$ws.Range("A91").Value = 'ankle support xs'
$ws.Range("A92").Value = 'anti cellulite jeans'
$ws.Range("A93").Value = 'anti see through leggings'
$ws.Range("A94").Value = 'antibacterial underwear men'
$ws.Range("A95").Value = 'antibacterial underwear women'
$ws.Range("A96").Value = 'apanx tights'
$ws.Range("A97").Value = 'apparel tape'
$ws.Range("A98").Value = 'apparel women'
$ws.Range("A99").Value = 'apperal for women'
$ws.Range("A100").Value = 'armour shorts'
